$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.169.25'
$ws.Range("E2").Value = '  +4.64%  '

$ws.Range("D3").Value = '4.070.13'
$ws.Range("E3").Value = '  +4.83%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '522.97'
$ws.Range("E5").Value = '  -1.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.79'
$ws.Range("E6").Value = '  +3.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.718'
$ws.Range("E7").Value = '  +18.06%  '

$ws.Range("D8").Value = '4.064.46'
$ws.Range("E8").Value = '  +4.89%  '

$ws.Range("E9").Value = '  +0.08%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.774'
$ws.Range("E10").Value = '  +7.99%  '

$ws.Range("E11").Value = '  +7.10%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000339'
$ws.Range("E12").Value = '  +3.31%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '48.77'
$ws.Range("E13").Value = '  +16.10%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.03'
$ws.Range("E14").Value = '  +7.34%  '

$ws.Range("D15").Value = '4.726.27'
$ws.Range("E15").Value = '  +5.35%  '

$ws.Range("D16").Value = '4.080.70'
$ws.Range("E16").Value = '  +6.08%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.45'
$ws.Range("E17").Value = '  +1.69%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '21.29'
$ws.Range("E18").Value = '  -1.28%  '

$ws.Range("E19").Value = '  +1.08%  '

$ws.Range("E20").Value = '  -0.04%  '

$ws.Range("D21").Value = '72.233.29'
$ws.Range("E21").Value = '  +4.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '448.95'
$ws.Range("E22").Value = '  +5.96%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '103.89'
$ws.Range("E23").Value = '  +18.85%  '

$ws.Range("E24").Value = '  +6.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.97'
$ws.Range("E25").Value = '  +5.90%  '

$ws.Range("E26").Value = '  +2.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.47'
$ws.Range("E27").Value = '  +1.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.16'
$ws.Range("E28").Value = '  +4.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '38.01'
$ws.Range("E29").Value = '  +5.14%  '

$ws.Range("E30").Value = '  +2.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.32'
$ws.Range("E31").Value = '  +16.58%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.71'
$ws.Range("E32").Value = '  +4.21%  '

$ws.Range("E33").Value = '  +5.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '687.19'
$ws.Range("E34").Value = '  -0.46%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '67.73'
$ws.Range("E35").Value = '  +1.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.66'
$ws.Range("E36").Value = '  +11.65%  '

$ws.Range("B37").Value = 'InjectiveProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '42.59'
$ws.Range("E37").Value = '  +6.25%  '

$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").Value = '0.0₃0903'
$ws.Range("E38").Value = '  +5.74%  '

$ws.Range("E39").Value = '  -1.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.154'
$ws.Range("E40").Value = '  +4.31%  '

$ws.Range("E41").Value = '  +7.80%  '

$ws.Range("E42").Value = '  -0.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0504'
$ws.Range("E43").Value = '  +4.96%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  +0.04%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.20'
$ws.Range("E45").Value = '  +0.24%  '

$ws.Range("E46").Value = '  +12.21%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.95'
$ws.Range("E47").Value = '  +16.25%  '

$ws.Range("E48").Value = '  -2.00%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000299'
$ws.Range("E49").Value = '  +8.73%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.40'
$ws.Range("E50").Value = '  +0.19%  '

$ws.Range("E51").Value = '  +4.05%  '
